$wb = $excel.ActiveWorkbook
$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

# --- Training Dashboard: updated "PERIOD TO EXPIRE" (H) values ---
$wsTraining.Range("H3").Value = 667
$wsTraining.Range("H4").Value = 378
$wsTraining.Range("H5").Value = 435
$wsTraining.Range("H6").Value = 482
$wsTraining.Range("H7").Value = 423
$wsTraining.Range("H8").Value = 594
$wsTraining.Range("H9").Value = 128
$wsTraining.Range("H10").Value = 244

# --- Training Dashboard: updated "LAST UPDATE" (I) dates, kept as plain text ---
$wsTraining.Range("I3:I10").NumberFormat = "@"
$wsTraining.Range("I3").Value = "16-Sep-2025"
$wsTraining.Range("I4").Value = "16-Sep-2025"
$wsTraining.Range("I5").Value = "16-Sep-2025"
$wsTraining.Range("I6").Value = "16-Sep-2025"
$wsTraining.Range("I7").Value = "16-Sep-2025"
$wsTraining.Range("I8").Value = "16-Sep-2025"
$wsTraining.Range("I9").Value = "16-Sep-2025"
$wsTraining.Range("I10").Value = "16-Sep-2025"

# --- Training Dashboard: header / title font now bold white ---
$wsTraining.Range("A1").Font.Bold = $true
$wsTraining.Range("A1").Font.Size = 11
$wsTraining.Range("A1").Font.Color = 16777215
$wsTraining.Range("A2:K2").Font.Bold = $true
$wsTraining.Range("A2:K2").Font.Color = 16777215

# --- Exam Dashboard: header / title font now bold white (same formatting table) ---
$wsExam.Range("A1").Font.Bold = $true
$wsExam.Range("A1").Font.Size = 11
$wsExam.Range("A1").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Bold = $true
$wsExam.Range("A2:G2").Font.Color = 16777215

# --- Exam Dashboard: widen CATEGORY/COMMENTS column and update remark text ---
# (14.1 is the ColumnWidth input that round-trips to a stored width of exactly 15
# through this host's char-width<->pixel rounding, matching the target column width)
$wsExam.Columns.Item(5).ColumnWidth = 14.1
$wsExam.Range("E3").Value = "date is valid"
